$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Highlight the "enrol" block (rows 12-14 and 16, column A:D) with the same
# green fill used elsewhere in the sheet (matches existing fill FF00B050).
# Note: COM Color uses BGR-in-int encoding => R=0,G=176,B=80 -> 0 + 176*256 + 80*65536
$greenBgr = 0 + (176 * 256) + (80 * 65536)

$ws.Range("A12:D14").Interior.Color = $greenBgr
$ws.Range("A16:D16").Interior.Color = $greenBgr

# Row 13 shrinks back down to the normal row height.
$ws.Rows("13:13").RowHeight = 17.25

# Remove the highlight (green fill + red font) from the "id crc error" row,
# turning it back into a plain/white row.
$ws.Range("A25:D25").Interior.ThemeColor = 2
$ws.Range("A25:D25").Font.ThemeColor = 1

# Move the active selection/cursor.
$ws.Range("C19").Select()
